$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1554434735375247
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 198602002.3250627
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 198602160.6238128
